# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Update the "last updated" timestamp text (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 11 de Abril de 2020 a las 09:52"

# --- Swap Irak / Estonia rows: Estonia now appears before Irak ---
# Row 65 was Irak, row 66 was Estonia. After the update, row 65 becomes
# Estonia (with refreshed data) and row 66 becomes Irak (with Irak's
# previous, unchanged data).
$ws.Range("A65").Value = "Estonia"
$ws.Range("A66").Value = "Irak"

# --- Row 29: Australia ---
$ws.Range("B29").Value = 6292
$ws.Range("C29").Value = 54
$ws.Range("D29").Value = 3265
$ws.Range("E29").Value = 2971
$ws.Range("F29").Value = 80
$ws.Range("G29").Value = 2
$ws.Range("H29").Value = 56

# --- Row 31: Polonia ---
$ws.Range("D31").Value = 375
$ws.Range("E31").Value = 5399

# --- Row 35: Rumania ---
$ws.Range("E35").Value = 4456
$ws.Range("G35").Value = 12
$ws.Range("H35").Value = 282

# --- Row 62: Moldavia ---
$ws.Range("D62").Value = 75
$ws.Range("E62").Value = 1334

# --- Row 65: now Estonia (refreshed data) ---
$ws.Range("B65").Value = 1304
$ws.Range("C65").Value = 46
$ws.Range("D65").Value = 93
$ws.Range("E65").Value = 1187
$ws.Range("F65").Value = 11
$ws.Range("G65").Value = 0
$ws.Range("H65").Value = 24

# --- Row 66: now Irak (carries the previous Irak data) ---
$ws.Range("B66").Value = 1279
$ws.Range("C66").Value = 0
$ws.Range("D66").Value = 550
$ws.Range("E66").Value = 659
$ws.Range("F66").Value = 0
$ws.Range("G66").Value = 0
$ws.Range("H66").Value = 70
